$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New StatQuery (Neo4j "count" query) text - replaces the old number_of_files
# query, now lives in column C for rows 2-5 (CasesTab/SamplesTab/FilesTab/
# StudyFilesTab all share the same StatQuery text).
# ---------------------------------------------------------------------------
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (f:file)-[*]->(samp:sample)-->(c)
MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp,demo, c, s, p, diag
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# ---------------------------------------------------------------------------
# New StudyFilesTab query (column B, row 5 - brand new tab).
# ---------------------------------------------------------------------------
$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (sf:file)-->(s)
MATCH (s)<--(c)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
WITH DISTINCT f,  s, c
WITH
        f, c,  s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c,  s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c,   s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# ---------------------------------------------------------------------------
# New CasesTab query (column B, row 2) - adds sample match, cohort, rounded
# age/weight and "Order by ... LIMIT 100" tail. Note the leading space and
# trailing spaces are part of the authored string.
# ---------------------------------------------------------------------------
$casesQuery = @'
 MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100        
'@

# ---------------------------------------------------------------------------
# Write order matters for shared-string table compaction/reindexing (mirrors
# how the authored workbook ended up with its particular string order): the
# StatQuery text first (frees+reuses the old stat-query slot across C2:C5),
# then the brand new row 5 (StudyFilesTab label + its query), and finally the
# new CasesTab query last, since it is the last thing to stop referencing the
# old (pre-edit) CasesTab text.
# ---------------------------------------------------------------------------

# Row 2 (CasesTab) StatQuery cell.
$ws.Range("C2").Value2 = $statQuery
$ws.Range("C2").WrapText = $true

# Row 3 (SamplesTab): query text unchanged, only the StatQuery in C changes.
$ws.Range("C3").Value2 = $statQuery
$ws.Range("C3").WrapText = $true

# Row 4 (FilesTab): query text unchanged, only the StatQuery in C changes.
$ws.Range("C4").Value2 = $statQuery
$ws.Range("C4").WrapText = $true

# Row 5 (new StudyFilesTab row).
$ws.Range("A5").Value2 = "StudyFilesTab"
$ws.Range("B5").Value2 = $studyFilesQuery
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value2 = $statQuery
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value2 = $ws.Range("D4").Value2
$ws.Range("E5").Value2 = $ws.Range("E4").Value2

# Row 2 (CasesTab) new query text (written last so it is the final new
# shared-string entry).
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("B2").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 285
$ws.Rows.Item(4).RowHeight = 270
$ws.Rows.Item(5).RowHeight = 409.5

# Row 6 (style-only row) - re-assert its wrap style.
$ws.Range("C6").WrapText = $true

# ---------------------------------------------------------------------------
# View: select C5 and scroll the window so row 5 area is visible (matches
# the updated tabSelected/topLeftCell/selection in the authored workbook).
# ---------------------------------------------------------------------------
$ws.Range("C5").Select()

# ---------------------------------------------------------------------------
# Column widths: nudge to the updated (slightly wider) bestfit widths; E
# widens substantially to fit the new StudyFilesTab content length.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.85546875
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 75.85546875
$ws.Columns.Item(4).ColumnWidth = 70.28515625
$ws.Columns.Item(5).ColumnWidth = 77.85546875
